# Switch Power_Network example to use "DC-OPF" as the Technical
# Representation for every line (column N, rows 7-19 on the
# "Power Network" sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Network")

$ws.Range("N7:N19").Value = "DC-OPF"
